# Apply the "handles float input without breaking stuff" fix to the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$xlCenter = -4108

function Set-Styled($addr, $styleName, $value) {
    $rng = $ws.Range($addr)
    $rng.Style = $styleName
    $rng.HorizontalAlignment = $xlCenter
    if ($null -ne $value) {
        $rng.Value = $value
    }
}

# --- Summary block (rows 10-12) ---------------------------------------
# Row 10: No. of Right / Wrong / Not-Attempt / Max
Set-Styled "A10" "mtitleStyle" "No."
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

# Row 11: Marking scheme (make sure the -1 is numeric, not text)
Set-Styled "A11" "mtitleStyle" "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Totals
Set-Styled "A12" "mtitleStyle" "Total"
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "53/112"

# --- Drop the third Student Ans / Correct Ans block (columns G:H) -----
$ws.Range("G1:H1048576").EntireColumn.Delete() | Out-Null

# --- Rework the per-question answer columns (A = result, D/E = 2nd set)
# Row 16 keeps a D/E pair, but it now reflects a wrong answer.
Set-Styled "A16" "correctStyle" "Option A"
Set-Styled "D16" "incorrectStyle" "Option B"

Set-Styled "A17" "correctStyle" "Option D"

Set-Styled "A18" "correctStyle" "Option B"

Set-Styled "A19" "correctStyle" "Option C"
$ws.Range("D19:E19").ClearContents()
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("D20:E20").ClearContents()
$ws.Range("D20:E20").Style = "Normal"

$ws.Range("D21:E21").ClearContents()
$ws.Range("D21:E21").Style = "Normal"

Set-Styled "A22" "correctStyle" "Option D"
$ws.Range("D22:E22").ClearContents()
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23:E23").ClearContents()
$ws.Range("D23:E23").Style = "Normal"

Set-Styled "A24" "correctStyle" "Option A"
$ws.Range("D24:E24").ClearContents()
$ws.Range("D24:E24").Style = "Normal"

Set-Styled "A25" "correctStyle" "Option A"
$ws.Range("D25:E25").ClearContents()
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("D26:E26").ClearContents()
$ws.Range("D26:E26").Style = "Normal"

Set-Styled "A27" "correctStyle" "Option A"
$ws.Range("D27:E27").ClearContents()
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D28:E28").ClearContents()
$ws.Range("D28:E28").Style = "Normal"

Set-Styled "A29" "correctStyle" "Option D"
$ws.Range("D29:E29").ClearContents()
$ws.Range("D29:E29").Style = "Normal"

Set-Styled "A30" "correctStyle" "Option B"
$ws.Range("D30:E30").ClearContents()
$ws.Range("D30:E30").Style = "Normal"

$ws.Range("D31:E31").ClearContents()
$ws.Range("D31:E31").Style = "Normal"

Set-Styled "A32" "correctStyle" "Option C"
$ws.Range("D32:E32").ClearContents()
$ws.Range("D32:E32").Style = "Normal"

Set-Styled "A33" "incorrectStyle" "Option A"
$ws.Range("D33:E33").ClearContents()
$ws.Range("D33:E33").Style = "Normal"

$ws.Range("D34:E34").ClearContents()
$ws.Range("D34:E34").Style = "Normal"

Set-Styled "A35" "correctStyle" "Option D"
$ws.Range("D35:E35").ClearContents()
$ws.Range("D35:E35").Style = "Normal"

Set-Styled "A36" "incorrectStyle" "Option D"
$ws.Range("D36:E36").ClearContents()
$ws.Range("D36:E36").Style = "Normal"

Set-Styled "A37" "correctStyle" "Option A"
$ws.Range("D37:E37").ClearContents()
$ws.Range("D37:E37").Style = "Normal"

$ws.Range("D38:E38").ClearContents()
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("D39:E39").ClearContents()
$ws.Range("D39:E39").Style = "Normal"

Set-Styled "A40" "correctStyle" "Option D"
$ws.Range("D40:E40").ClearContents()
$ws.Range("D40:E40").Style = "Normal"
